$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.926.01"
$ws.Range("E2").Value = "  +2.25%  "
$ws.Range("D3").Value = "'3.593.21"
$ws.Range("E3").Value = "  +1.49%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'606.41"
$ws.Range("E5").Value = "  +4.35%  "
$ws.Range("D6").Value = "'175.13"
$ws.Range("E6").Value = "  +1.39%  "
$ws.Range("D7").Value = "'3.590.05"
$ws.Range("E7").Value = "  +1.60%  "
$ws.Range("E8").Value = "  +1.23%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("E10").Value = "  +4.57%  "
$ws.Range("D11").Value = "'7.43"
$ws.Range("E11").Value = "  +9.08%  "
$ws.Range("D12").Value = "'0.592"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "'47.33"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("D14").Value = "'0.0000279"
$ws.Range("E14").Value = "  +1.33%  "
$ws.Range("D15").Value = "'4.171.67"
$ws.Range("E15").Value = "  +1.52%  "
$ws.Range("D16").Value = "'8.48"
$ws.Range("E16").Value = "  -0.68%  "
$ws.Range("D17").Value = "'622.71"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "'3.629.79"
$ws.Range("E18").Value = "  +2.30%  "
$ws.Range("D19").Value = "'71.174.74"
$ws.Range("E19").Value = "  +2.53%  "
$ws.Range("E20").Value = "  -2.59%  "
$ws.Range("D21").Value = "'17.53"
$ws.Range("E21").Value = "  +0.73%  "
$ws.Range("D22").Value = "'0.892"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'9.39"
$ws.Range("E23").Value = "  -16.14%  "
$ws.Range("D24").Value = "'16.21"
$ws.Range("E24").Value = "  +1.52%  "
$ws.Range("D25").Value = "'98.04"
$ws.Range("E25").Value = "  +0.76%  "
$ws.Range("D26").Value = "'3.82"
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  +0.98%  "
$ws.Range("D29").Value = "'9.36"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").Value = "'33.76"
$ws.Range("E30").Value = "  +2.63%  "
$ws.Range("D31").Value = "'8.56"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("D32").Value = "'3.09"
$ws.Range("E32").Value = "  -2.23%  "
$ws.Range("D33").Value = "'7.15"
$ws.Range("E33").Value = "  +2.18%  "
$ws.Range("D34").Value = "'1.32"
$ws.Range("E34").Value = "  -2.07%  "
$ws.Range("D35").Value = "'627.74"
$ws.Range("E35").Value = "  -1.84%  "
$ws.Range("D36").Value = "'3.79"
$ws.Range("E36").Value = "  +7.69%  "
$ws.Range("E37").Value = "  -0.25%  "
$ws.Range("D38").Value = "'10.91"
$ws.Range("E38").Value = "  +1.00%  "
$ws.Range("D39").Value = "'0.0485"
$ws.Range("E39").Value = "  +6.44%  "
$ws.Range("D40").Value = "'57.44"
$ws.Range("E40").Value = "  +0.29%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.25%  "
$ws.Range("E42").Value = "  +4.49%  "
$ws.Range("D43").Value = "'3.421.02"
$ws.Range("E43").Value = "  +0.92%  "
$ws.Range("E44").Value = "  -0.75%  "
$ws.Range("D45").Value = "'0.0₃0720"
$ws.Range("E45").Value = "  +2.56%  "
$ws.Range("E46").Value = "  +9.08%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.71"
$ws.Range("E47").Value = "  +5.97%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'33.22"
$ws.Range("E48").Value = "  +1.32%  "
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("D50").Value = "'132.65"
$ws.Range("E50").Value = "  +0.09%  "

